$d = $word.ActiveDocument

# The document has two paragraphs whose runs spell out
#   "你好，" + "$\SunQuarTeX$" [+ "!"]
# split across 2 (resp. 3) separate <w:r> runs. Collapse each such paragraph
# down to a single run whose text has the literal "$...$" math-delimiters
# around SunQuarTeX removed, e.g. "你好，SunQuarTeX" / "你好，SunQuarTeX!".

$needle = '$\SunQuarTeX$'
$replacement = 'SunQuarTeX'

foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text
    if ($paraText.Contains($needle)) {
        # Paragraph.Range.Text includes the trailing paragraph-mark char;
        # strip it off so we only rewrite the visible text.
        $body = $paraText.Substring(0, $paraText.Length - 1)

        $newBody = $body.Replace($needle, $replacement)

        # Re-target a range over just the paragraph's text (not its ending
        # mark) and overwrite it; Word collapses this into a single run
        # that inherits the formatting of the paragraph's first run.
        $bodyRange = $d.Range($p.Range.Start, $p.Range.End - 1)
        $bodyRange.Text = $newBody
    }
}
